# Apply updated cryptocurrency price/volume data to worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D holds text values such as "1.001" or "28.011.06" that must not
# be auto-converted to numbers by Excel, so force the column to Text format
# before writing any of the new values into it.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range('D2').Value = '28.011.06'
$ws.Range('E2').Value = '  -1.96%  '
$ws.Range('D3').Value = '1.830.05'
$ws.Range('D4').Value = '1.001'
$ws.Range('E4').Value = '  -0.21%  '
$ws.Range('D5').Value = '324.74'
$ws.Range('E5').Value = '  -3.21%  '
$ws.Range('E6').Value = '  -0.19%  '
$ws.Range('E7').Value = '  -0.46%  '
$ws.Range('D8').Value = '0.3874'
$ws.Range('E8').Value = '  -1.13%  '
$ws.Range('D9').Value = '0.07862'
$ws.Range('E9').Value = '  -0.28%  '
$ws.Range('D10').Value = '0.9587'
$ws.Range('E10').Value = '  -2.66%  '
$ws.Range('D11').Value = '21.82'
$ws.Range('E11').Value = '  -1.65%  '
$ws.Range('D12').Value = '1.835.63'
$ws.Range('E12').Value = '  -3.10%  '
$ws.Range('D13').Value = '5.671'
$ws.Range('E13').Value = '  -3.10%  '
$ws.Range('D14').Value = '6.899'
$ws.Range('E14').Value = '  -1.78%  '
$ws.Range('D15').Value = '0.06785'
$ws.Range('E15').Value = '  -0.66%  '
$ws.Range('D16').Value = '87.23'
$ws.Range('E16').Value = '  -0.49%  '
$ws.Range('D17').Value = '1.002'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').Value = '0.000009909'
$ws.Range('E18').Value = '  -1.93%  '
$ws.Range('D19').Value = '16.61'
$ws.Range('E19').Value = '  -2.47%  '
$ws.Range('D20').Value = '1.001'
$ws.Range('E20').Value = '  -0.06%  '
$ws.Range('D21').Value = '28.017.73'
$ws.Range('E21').Value = '  -2.00%  '
$ws.Range('D22').Value = '5.314'
$ws.Range('E22').Value = '  -1.77%  '
$ws.Range('D23').Value = '10.97'
$ws.Range('E23').Value = '  -2.63%  '
$ws.Range('D24').Value = '2.092'
$ws.Range('E24').Value = '  -1.43%  '
$ws.Range('D25').Value = '2.070.03'
$ws.Range('E25').Value = '  -2.37%  '
$ws.Range('D26').Value = '153.80'
$ws.Range('E26').Value = '  +0.13%  '
$ws.Range('D27').Value = '19.15'
$ws.Range('E27').Value = '  -1.14%  '
$ws.Range('D28').Value = '5.734'
$ws.Range('E28').Value = '  -7.01%  '
$ws.Range('D29').Value = '1.972'
$ws.Range('E29').Value = '  -2.51%  '
$ws.Range('D30').Value = '117.42'
$ws.Range('E30').Value = '  -0.24%  '
$ws.Range('B31').Value = 'Stellar'
$ws.Range('C31').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D31').Value = '0.09258'
$ws.Range('E31').Value = '  -1.96%  '
$ws.Range('B32').Value = 'ImmutableX'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D32').Value = '0.9352'
$ws.Range('E32').Value = '  -4.19%  '
$ws.Range('D33').Value = '5.290'
$ws.Range('E33').Value = '  -1.62%  '
$ws.Range('E34').Value = '  -2.29%  '
$ws.Range('D35').Value = '3.286'
$ws.Range('E35').Value = '  -6.16%  '
$ws.Range('D36').Value = '0.05864'
$ws.Range('E36').Value = '  -4.12%  '
$ws.Range('D37').Value = '0.02142'
$ws.Range('E37').Value = '  -2.46%  '
$ws.Range('D38').Value = '1.146'
$ws.Range('E38').Value = '  -1.41%  '
$ws.Range('D39').Value = '7.767'
$ws.Range('E39').Value = '  +2.35%  '
$ws.Range('D40').Value = '0.5580'
$ws.Range('E40').Value = '  -2.19%  '
$ws.Range('D41').Value = '9.870'
$ws.Range('E41').Value = '  -2.47%  '
$ws.Range('D42').Value = '0.1762'
$ws.Range('E42').Value = '  -1.68%  '
$ws.Range('D43').Value = '11.60'
$ws.Range('E43').Value = '  -2.07%  '
$ws.Range('D44').Value = '0.5262'
$ws.Range('E44').Value = '  -2.42%  '
$ws.Range('D45').Value = '0.07001'
$ws.Range('E45').Value = '  -2.21%  '
$ws.Range('D46').Value = '2.122'
$ws.Range('E46').Value = '  -10.74%  '
$ws.Range('D47').Value = '1.115'
$ws.Range('E47').Value = '  -11.56%  '
$ws.Range('D48').Value = '1.829'
$ws.Range('E48').Value = '  -4.12%  '
$ws.Range('D49').Value = '112.97'
$ws.Range('E49').Value = '  -0.73%  '
$ws.Range('D50').Value = '1.000'
$ws.Range('E50').Value = '  -0.14%  '
$ws.Range('D51').Value = '2.319'
$ws.Range('E51').Value = '  +0.24%  '
